# The sheet currently has an empty "header" row (row 6) that only carries
# the label "grandes regiões e unidades da federação" in column A, with no
# data in B:G. The fix removes that stray label/row entirely: every row
# below it (the actual region/state data rows 7-38) shifts up by one, the
# now-unused shared string is dropped automatically on save, and the
# worksheet's used range shrinks from G38 to G37.
#
# That is exactly what a native "delete entire row" does, so we just
# delete row 6 and let Excel shift everything below it upward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Delete()
